# Re-sort the needle calibration data (rows 2-8) in ascending order of
# column A (time), as performed during calibration of the needle.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @([double]"52899.342991", [double]"-3.942213424e-05",  [double]"-3.0275960474e-05", [double]"-3.2085221738e-06"),
    @([double]"52922.142992", [double]"-0.00020350871801", [double]"-0.0001621936335",   [double]"-3.9442408426e-05"),
    @([double]"52933.142993", [double]"-0.00046556799422", [double]"-0.00037468526522",  [double]"-9.4001703476e-05"),
    @([double]"52943.674993", [double]"-0.0007301185",     [double]"-0.0005968197",      [double]"-0.0001474277"),
    @([double]"52953.142994", [double]"-0.0004484188",     [double]"-0.0003643707",      [double]"-0.000104648"),
    @([double]"52963.406994", [double]"-0.00018641813817", [double]"-0.00015060479696",  [double]"-4.4232515122e-05"),
    @([double]"52973.938995", [double]"-3.2965453674e-05", [double]"-2.6705059835e-05",  [double]"-8.1586853441e-06")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}
